$wb = $excel.ActiveWorkbook

# --- Sheet "record": only the active cell/selection changes ---
$wsRecord = $wb.Worksheets.Item("record")
$wsRecord.Activate() | Out-Null
$wsRecord.Range("A1210").Select() | Out-Null

# --- Sheet "books": scroll position + active cell changes ---
$wsBooks = $wb.Worksheets.Item("books")
$wsBooks.Activate() | Out-Null
$winBooks = $excel.ActiveWindow
$winBooks.ScrollRow = 1
$winBooks.ScrollColumn = 1
$wsBooks.Range("B6").Select() | Out-Null

# --- Sheet "names": rebuild as a single sorted column of unique names ---
$wsNames = $wb.Worksheets.Item("names")
$wsNames.Activate() | Out-Null

$names = @(
    "Aazadh",
    "Abdul rahoof s",
    "Abhinand c. S.",
    "Abhinaya",
    "Ajay pradeep",
    "Akhil vasim",
    "Akhila s",
    "Akhiya a.r.",
    "Alan syed",
    "Amritha raju",
    "Amrutha p",
    "Anagha rajan",
    "Anas a s",
    "Ann",
    "Anu sasi",
    "Anuja mohan",
    "Anuroopa g nadh",
    "Arjun v.m.",
    "Arun k unni",
    "Arun paulose",
    "Arunima",
    "Arunkumar b",
    "Arya k r",
    "Arya krishna a r",
    "Arya v. V.",
    "Asha v s",
    "Ashlin",
    "Ashna",
    "Ashwini jayachandran",
    "Ashwini k.s.",
    "Aswathy t r",
    "Athira h",
    "Athira t",
    "Banazir",
    "Bibin",
    "Bishmitha k",
    "Chinchu e.r.",
    "Christy",
    "Deepthi",
    "Delphin r. D.",
    "Dineshkumar",
    "Dr. Achuthsankar s. Nair",
    "Dr. Oommen v oommen",
    "Gayathri d",
    "Gayathri thampi",
    "Gopika rani",
    "Gowri k babu",
    "Jijil",
    "Junaida m i",
    "Karthika",
    "Krishnendu",
    "Lekshmi s nair",
    "Lidhiya",
    "Liya simon",
    "Madeena s",
    "Manumol m",
    "Muhammed hakkim",
    "Muhammed minhaju a.",
    "Neema m p",
    "Neenu mohan",
    "Nighitha",
    "Parvathi jayaraj",
    "Parvathy suresh",
    "Preetha p.",
    "Rani j.r.",
    "Raniya k zubair",
    "Rashmi sukumaran",
    "Roshan",
    "Saleena younus",
    "Salil suresh",
    "Sameera k.",
    "Sandhya k. S.",
    "Saranya s",
    "Saraswathy v",
    "Sarath kumar",
    "Seba",
    "Shahana s j",
    "Shahina k",
    "Shajila salim",
    "Shanitha a.",
    "Sheeba k",
    "Siva priya p.",
    "Sneha thomas",
    "Sreejith g",
    "Sreelekshmi i.g.",
    "Sunitha p.",
    "Suveena s",
    "Swathi k",
    "Vaishnavi",
    "Vijayalekshmi b",
    "Vinni n g",
    "Vinod m.p.",
    "Vishnu",
    "Vishnu v j"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $wsNames.Cells.Item($i + 1, 1).Value = $names[$i]
}

# Remove the now-unused rows 95:96 (old sheet had 96 rows, new has 94)
$wsNames.Range("A95:B96").EntireRow.Delete()

# Remove column B entirely (values moved into column A only)
$wsNames.Columns.Item(2).Delete()

# Match the new column width for column A
$wsNames.Columns.Item(1).ColumnWidth = 24.866666666666667

$winNames = $excel.ActiveWindow
$winNames.ScrollRow = 60
$winNames.ScrollColumn = 1
$wsNames.Range("A100").Select()

# Re-activate the originally active sheet ("record") last, so it remains
# the selected/active tab in the saved workbook.
$wsRecord.Activate()
